# Commit: "added link to repo in pptx"
#
# 1. Insert a new slide at position 2 (Title and Content layout) with a
#    short "3 Things:" slide pointing at the repo / W3Schools / MDN.
# 2. Tidy up slide 1's subtitle so it is a single run instead of three.

$p = $ppt.ActivePresentation

# --- 1. New slide at index 2, "Title and Content" layout (ppLayoutText = 2) ---
$ns = $p.Slides.Add(2, 2)

$ns.Shapes.Item(1).TextFrame.TextRange.Text = "3 Things:"

$body = $ns.Shapes.Item(2).TextFrame
$body.AutoSize = 0   # ppAutoSizeNone -> <a:noAutofit/>

$tr = $body.TextRange
$tr.Text = "The repository is here: "
$runHttp = $tr.InsertAfter("http://")
$runLink = $runHttp.InsertAfter("bit.ly/2mRzSc3")
$runLink.InsertAfter("`r")
$runLink.InsertAfter("W3 Schools " + [char]0x2013 + " a great resource for HTML/CSS and a little JS")
$tr.Paragraphs(2).InsertAfter("`r")
$tr.Paragraphs(3).InsertAfter("MDN (Mozilla Developer Network) " + [char]0x2013 + " Best ECMAScript/JavaScript reference out there")

# Sizes / bold for paragraph 1 (url bold, rest of body text sz 40)
$full = $body.TextRange
$full.Font.Size = 40
$lenIntro = "The repository is here: ".Length
$lenHttp = "http://".Length
$lenLink = "bit.ly/2mRzSc3".Length
$full.Characters($lenIntro + 1, $lenHttp + $lenLink).Font.Bold = $true

# --- 2. Slide 1 subtitle: merge the three runs into one ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitle.Text = "placeholder"
$subtitle2 = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitle2.Text = "An introductory workshop on the basic components of a webpage"
